# Applies text replacements to the "two-digit number multiplied by
# two-digit number" worksheet document.

$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-30 Tuesday", "2025-12-31 Wednesday"),
    @("24×83=", "22×43="),
    @("40×71=", "93×30="),
    @("82×73=", "20×23="),
    @("79×84=", "99×17="),
    @("31×47=", "52×23="),
    @("13×17=", "37×41="),
    @("15×18=", "40×52="),
    @("13×94=", "82×48="),
    @("91×70=", "92×15="),
    @("65×17=", "15×16="),
    @("72×61=", "49×82="),
    @("67×61=", "57×93="),
    @("29×97=", "97×89="),
    @("63×60=", "54×39="),
    @("17×80=", "36×20="),
    @("51×73=", "69×70="),
    @("91×84=", "63×17="),
    @("59×42=", "54×73="),
    @("35×17=", "35×76="),
    @("88×35=", "24×88="),
    @("45×53=", "64×88="),
    @("91×30=", "71×97="),
    @("70×58=", "54×18="),
    @("39×81=", "24×78="),
    @("88×57=", "84×81=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
